$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New S-val data (regenerated to filter save games)
$data = @{
    2 = @{ B = 3.182878228561681;  C = 1.65323645889881;  D = 0.7127328510149897;  E = 6.48142807727062;   G = 12.0302756157461 }
    3 = @{ B = 3.182878228561681;  C = 1.65323645889881;  D = 0.1529057820181812;  E = 0.4998867070740569; G = 5.488907176552729 }
    4 = @{ B = 0.1554434735375247; C = 0.3375848360084654; D = 16.98373111632243;  E = 0.4998867070740569; G = 17.97664613294248 }
    5 = @{ B = 0.3464964993005633; C = 0.3375848360084654; D = 3.082599426703578;  E = 0.4998867070740569; G = 4.266567469086664 }
    6 = @{ B = 1.505614041169197;  C = 1.65323645889881;  D = 0.1529057820181812;  E = 0.4998867070740569; G = 3.811642989160245 }
    7 = @{ B = 3.182878228561681;  C = 1.65323645889881;  D = 0.7127328510149897;  E = 0.4998867070740569; G = 6.048734245549538 }
    8 = @{ B = 1.505614041169197;  C = 1.65323645889881;  D = 0.1529057820181812;  E = 0.4998867070740569; G = 3.811642989160245 }
    9 = @{ B = 0.06328177979961902; C = 0.3375848360084654; D = 0.7127328510149897; E = 0.4998867070740569; G = 1.613486173897131 }
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Cells.Item($row, 2).Value = $vals.B
    $ws.Cells.Item($row, 3).Value = $vals.C
    $ws.Cells.Item($row, 4).Value = $vals.D
    $ws.Cells.Item($row, 5).Value = $vals.E
    $ws.Cells.Item($row, 7).Value = $vals.G
}
